# Excel to txt seems to be working.
#
# The source data only added one new shared string ("Lorem ...") that gets
# written into the new column B for every data row, and the saved selection
# moved from E1 to B3 (presumably where the cursor ended up after typing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Populate the new column B (rows 2-4) with the new value.
$ws.Range("B2:B4").Value = "Lorem ..."

# Leave the selection where Excel would have left it after entering the data.
$ws.Range("B3").Select()
